# Applies the "456a3b4" content refresh to the 展览 (Exhibitions) and
# 全部类型 (All types) sheets: a handful of "want-to-go" counter bumps,
# plus three newly scraped events spliced into the May/July block
# (which pushes every following row down).

function Set-IndexStyle($ws, $row) {
    # Column A carries a running index (header=0, then 1,2,3,...) and uses
    # the workbook's bold/centered/bordered "index" style. Row-Insert()
    # clones formatting from the row above/below but lands on a *new*
    # style slot; re-asserting these properties collapses it back onto
    # the single shared style every other index cell already uses.
    $c = $ws.Cells.Item($row, 1)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Borders.LineStyle = 1         # xlContinuous
}

function Set-EventRow($ws, $row, $b, $c, $d, $e, $f, $g, $h, $i) {
    Set-IndexStyle $ws $row
    $ws.Cells.Item($row, 1).Value = $row - 1
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

function Update-ExhibitionSheet($ws, $rowOffset) {
    # --- plain "want-to-go" counter refreshes, earlier in the sheet ---
    $ws.Cells.Item(3 + $rowOffset, 6).Value = 2225
    $ws.Cells.Item(4 + $rowOffset, 6).Value = 96
    $ws.Cells.Item(5 + $rowOffset, 6).Value = 13347
    if ($rowOffset -eq 0) {
        $ws.Cells.Item(9, 6).Value = 485
        $ws.Cells.Item(10, 6).Value = 1196
        $ws.Cells.Item(11, 6).Value = 1004
        $ws.Cells.Item(12, 6).Value = 13820
        $ws.Cells.Item(13, 6).Value = 14480
    } else {
        # "全部类型" has one extra row (a 演出 item) inserted above this
        # block already, so these five sit one row further down.
        $ws.Cells.Item(10 + $rowOffset, 6).Value = 485
        $ws.Cells.Item(11 + $rowOffset, 6).Value = 1196
        $ws.Cells.Item(12 + $rowOffset, 6).Value = 1004
        $ws.Cells.Item(13 + $rowOffset, 6).Value = 13820
        $ws.Cells.Item(14 + $rowOffset, 6).Value = 14480
    }

    # Base row number (in the "展览" sheet) where the new cv meet-and-greet
    # event is inserted; add $rowOffset for "全部类型".
    $r1 = 21 + $rowOffset

    # 1) New row: 苏州·动漫游戏嘉年华cv见面会
    $ws.Rows.Item($r1).Insert()
    Set-EventRow $ws $r1 "2024-05-02" "苏州·动漫游戏嘉年华cv见面会" `
        "东太湖大道12000号 苏州湾大剧院" "2024.05.02 14:00-05.02 17:00" 2 168 `
        "https://show.bilibili.com/platform/detail.html?id=83504" `
        "//i2.hdslb.com/bfs/openplatform/202403/YqMudwaj1711608967902.jpeg"

    # 2) New row: 太仓·龙狮动漫嘉年华5.0 (goes in two rows further down,
    #    after the untouched "苏州·苏州湾动漫游戏嘉年华" row)
    $r2 = $r1 + 2
    $ws.Rows.Item($r2).Insert()
    Set-EventRow $ws $r2 "2024-05-03" "太仓·龙狮动漫嘉年华5.0" `
        "滨河路128号 凯景世纪大酒店(太仓滨河路店)" "2024.05.03 08:00-05.03 17:00" 1 45 `
        "https://show.bilibili.com/platform/detail.html?id=83507" `
        "//i2.hdslb.com/bfs/openplatform/202403/reGF2YIi1711420063540.png"

    # Rows between the second new row and the third new row keep their old
    # content but several "want-to-go" counters ticked up; row numbers below
    # already include both inserts above.
    Set-EventRow $ws ($r2 + 2) "2024-05-03" "昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会" `
        "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.03 14:00-05.03 16:00" 116 1 `
        "https://show.bilibili.com/platform/detail.html?id=81120" `
        "//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg"

    Set-EventRow $ws ($r2 + 3) "2024-05-03" "昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会" `
        "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.03 14:00-05.03 16:00" 58 1 `
        "https://show.bilibili.com/platform/detail.html?id=81114" `
        "//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg"

    Set-EventRow $ws ($r2 + 4) "2024-05-04" "苏州·OCG国潮动漫游戏嘉年华" `
        "苏州大道东688号 苏州国际博览中心" "2024.05.04 09:00-05.05 17:00" 5533 65 `
        "https://show.bilibili.com/platform/detail.html?id=82779" `
        "//i1.hdslb.com/bfs/openplatform/202403/hcgdIzw61710298907237.jpeg"

    Set-EventRow $ws ($r2 + 5) "2024-05-04" "苏州·OCG国潮动漫游戏嘉年华阿杰内场" `
        "苏州大道东688号 苏州国际博览中心" "2024.05.04 09:00-05.04 17:00" 943 "已售罄" `
        "https://show.bilibili.com/platform/detail.html?id=82940" `
        "//i2.hdslb.com/bfs/openplatform/202403/lLKmv48C1710511298160.jpeg"

    Set-EventRow $ws ($r2 + 6) "2024-05-18" "苏州·YoungComic动漫嘉年华" `
        "清禾路886号 尹山湖大剧院" "2024.05.18 10:00-05.18 17:00" 1036 60 `
        "https://show.bilibili.com/platform/detail.html?id=83142" `
        "//i2.hdslb.com/bfs/openplatform/202403/4wWLK6Jg1710840463319.jpeg"

    Set-EventRow $ws ($r2 + 7) "2024-06-08" "【会员购严选】苏州·Come in joy动漫国潮文化节" `
        "金山南路288号 广电国际会展中心" "2024.06.08 10:00-06.09 17:00" 2332 60 `
        "https://show.bilibili.com/platform/detail.html?id=82233" `
        "//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg"

    Set-EventRow $ws ($r2 + 8) "2024-06-29" "苏州·归离之缘原神only展" `
        "清禾路888号2号楼3楼 格莱美婚礼宴会中心" "2024.06.29 09:30-06.29 18:30" 30 89 `
        "https://show.bilibili.com/platform/detail.html?id=83271" `
        "//i1.hdslb.com/bfs/openplatform/202403/hNkSoRCt1710999968954.png"

    # 3) New row: 苏州·白日梦想7.20全职猎人ONLY展, inserted right before the
    #    trailing 萤火国潮文化节 row.
    $r3 = $r2 + 9
    $ws.Rows.Item($r3).Insert()
    Set-EventRow $ws $r3 "2024-07-20" "苏州·白日梦想7.20全职猎人ONLY展" `
        "金芳路与新发路交叉口东南120米 万龙大厦" "2024.07.20 09:00-07.20 17:00" 14 72 `
        "https://show.bilibili.com/platform/detail.html?id=83508" `
        "//i1.hdslb.com/bfs/openplatform/202403/LC3LtiWw1711633827120.jpeg"

    Set-EventRow $ws ($r3 + 1) "2024-07-20" "苏州·萤火国潮文化节动漫品牌博览会" `
        "金山南路288号木渎影视城F2 苏州广电国际会展中心" "2024.07.20 10:00-07.21 17:00" 117 60 `
        "https://show.bilibili.com/platform/detail.html?id=83301" `
        "//i0.hdslb.com/bfs/openplatform/202403/rV07luU61711274774556.jpeg"
}

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
Update-ExhibitionSheet $wsExhibitions 0

$wsAll = $wb.Worksheets.Item("全部类型")
Update-ExhibitionSheet $wsAll 1
